$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProjectSchedule")
Write-Host "B1 before: $($ws.Range("B1").Value)"
$ws.Range("B1").Value = "E-commerce Website "
Write-Host "B1 after: $($ws.Range("B1").Value)"
